# "initial level select touch-up"
# Insert a new "okay" / level-match-title / level-match-desc localization
# block right after the "close" row (row 11), pushing everything below it
# down by three rows, then fill in the new strings and move the on-sheet
# selection to line up with where the author was working (B15, the new
# "cycle" row that used to be B12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before row 12 (old row 12 "cycle" shifts to row 15).
$ws.Rows("12:14").Insert(-4121, 0)

# New localization keys/values, column A then column B per row (matches the
# order new shared strings were appended in the source edit).
$ws.Range("A12").Value = "okay"
$ws.Range("B12").Value = "OKAY"
$ws.Range("A13").Value = "levelMatchTitle"
$ws.Range("A14").Value = "levelMatchDesc"
$ws.Range("B13").Value = "MATCH CLIMATE"
$ws.Range("B14").Value = "Find the location that matches the following climate zone and type."

# Scroll back to the top and leave the selection where the author left it.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("B15").Select() | Out-Null

Write-Host "done"
